# Adds fiscal-year rows 2021 and 2022 (sheet rows 11-12) to the "管理费用"
# (management expense) time series, mirroring the existing row layout:
#   column A  -> "<year>年" label, styled like the other year cells
#   columns B:AQ -> one management-expense figure per industry column
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- row 11: 2021年 ----
# Clone the formatting of the previous year-label cell (bold / border / centered)
$ws.Range("A10").Copy($ws.Range("A11"))
$ws.Range("A11").Value = "2021年"
$row11 = @(3357.28, 1045.06, 213.61, 1.43, 1635.32, 4363.09, 401.85, 3051.24, 564.39, 564.99, 226.73, 118.52, 809.4, 1688.5, 286.03, 379.85, 1970.38, 399.07, 5441.97, 849.97, 2104.72, 369.89, 1742.5, 5314.95, 608.99, 1255.43, 770.55, 1274.46, 897.37, 69764.37, 10043.85, 3831.04, 787.75, 916.92, 162.22, 2612.92, 1476.78, 3586.34, 312, 1176.83, 2791.44, 358.72)
for ($i = 0; $i -lt $row11.Length; $i++) {
    $ws.Cells.Item(11, $i + 2).Value = $row11[$i]
}

# ---- row 12: 2022年 ----
# Clone the formatting of the previous year-label cell (bold / border / centered)
$ws.Range("A10").Copy($ws.Range("A12"))
$ws.Range("A12").Value = "2022年"
$row12 = @(3520, 1071.6, 148.9, 1.1, 1576.6, 4661.2, 422.1, 3264, 567.8, 561.7, 237.5, 121.8, 805.1, 1710.9, 290, 373.7, 1984.9, 413.5, 5443.2, 825.9, 2387.7, 394.2, 1788.4, 5808, 605.9, 1176.7, 795.1, 1304.4, 870.9, 71238.8, 10607, 3871.3, 799, 916.6, 166.9, 2593.5, 1082.9, 3625.4, 319.6, 1237.1, 2594.3, 292.5)
for ($i = 0; $i -lt $row12.Length; $i++) {
    $ws.Cells.Item(12, $i + 2).Value = $row12[$i]
}
